$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55 currently holds phone "09876543" stored as text (with a leading
# zero). The update normalizes A55 to the numeric phone 9876543, and the
# original text phone "09876543" is preserved as a new row 56 (same
# birthday/points), growing the sheet from A1:C55 to A1:C56.

# Row 55: phone becomes a plain number (loses the leading zero), points stay 0.
$ws.Cells.Item(55, 1).Value = 9876543
$ws.Cells.Item(55, 3).Value = 0

# Row 56: new row - recreate the original text phone, blank birthday, 0 points.
$ws.Cells.Item(56, 1).Value = "'09876543"
$ws.Cells.Item(56, 1).ClearFormats()
$ws.Cells.Item(56, 2).Value = "'"
$ws.Cells.Item(56, 2).ClearFormats()
$ws.Cells.Item(56, 3).Value = 0
